$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 104: TravelBackTooltip
$ws.Range("A104").Value = "TravelBackTooltip"
$ws.Range("B104").Value = "The saw [00FFFF]returns[-] to its original position by [FFFF00]passing through the waypoints.[-]"
$ws.Range("C104").Value = "La sierra [00FFFF]regresa[-] a su posición original al [FFFF00]pasar por todos sus waypoints.[-]"
$ws.Rows.Item(104).RowHeight = 30

# Row 105: LoopTooltip
$ws.Range("A105").Value = "LoopTooltip"
$ws.Range("B105").Value = "The saw [00FFFF]returns[-] to its original position in a [FFFF00]straight line [b]ignoring[/b] the waypoints.[-]"
$ws.Range("C105").Value = "La sierra [00FFFF]regresa[-] a su posición original en una [FFFF00]linea recta [b]ignorando[/b] los waypoints.[-]"
$ws.Rows.Item(105).RowHeight = 45

$ws.Range("B105").Select() | Out-Null
